# Commit message: "added check for correct input in add_worker_Constraints"
# This adds two new per-worker constraint sheets ("tair" and "asaf"), following
# the same layout already used for the existing "yoni" constraints sheet:
#   Row1: B1:H1 = Sunday..Saturday (day-of-week header)
#   Row2: A2 = "Morning", with "NO" placed under the day(s) the worker can't work mornings
#   Row3: A3 = "Evening", with "NO" placed under the day(s) the worker can't work evenings

$wb = $excel.ActiveWorkbook

function Add-ConstraintsSheet($wb, $sheetName, $morningNoCols, $eveningNoCols) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $sheetName

    $ws.Range("B1").Value = "Sunday"
    $ws.Range("C1").Value = "Monday"
    $ws.Range("D1").Value = "Tuesday"
    $ws.Range("E1").Value = "Wednesday"
    $ws.Range("F1").Value = "Thursday"
    $ws.Range("G1").Value = "Friday"
    $ws.Range("H1").Value = "Saturday"

    $ws.Range("A2").Value = "Morning"
    foreach ($col in $morningNoCols) {
        $ws.Range($col + "2").Value = "NO"
    }

    $ws.Range("A3").Value = "Evening"
    foreach ($col in $eveningNoCols) {
        $ws.Range($col + "3").Value = "NO"
    }

    return $ws
}

# Sheet "tair": morning NO on Monday (C2); evening NO on Sunday (B3)
Add-ConstraintsSheet $wb "tair" @("C") @("B") | Out-Null

# Sheet "asaf": morning NO on Monday (C2) and Friday (G2); no evening NO values
Add-ConstraintsSheet $wb "asaf" @("C", "G") @() | Out-Null
